$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Color constants (OLE BGR longs)
# -----------------------------------------------------------------
$cBlack      = 0         # FF000000
$cFillHeader = 11711407  # FFAFB3B2
$cBorderLt   = 9737364   # FF949494
$cBorderDk   = 3158064   # FF303030
$cBorderLt2  = 10132122  # FF9A9A9A

$xlEdgeLeft   = 7
$xlEdgeTop    = 8
$xlEdgeBottom = 9
$xlEdgeRight  = 10
$xlContinuous = 1

# -----------------------------------------------------------------
# 1) Add the new "isVerified" column (H) with template placeholders
# -----------------------------------------------------------------
$ws.Range("H1").Value = "{d.i18n.isVerified}"
$ws.Range("H2").Value = "{d.contacts[i].isVerified}"
$ws.Range("H3").Value = "{d.contacts[i+1].isVerified}"

# -----------------------------------------------------------------
# 2) Header row (row 1): bold Arial/10/black, solid grey fill,
#    thin borders - light on L/T/R, dark on bottom
# -----------------------------------------------------------------
for ($col = 1; $col -le 8; $col++) {
  $cell = $ws.Cells.Item(1, $col)
  $cell.Font.Name = "Arial"
  $cell.Font.Size = 10
  $cell.Font.Bold = $true
  $cell.Font.Color = $cBlack
  $cell.Interior.Pattern = 1
  $cell.Interior.Color = $cFillHeader
  $cell.NumberFormat = "General"
  $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeLeft).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeTop).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeRight).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeBottom).Color = $cBorderDk
}

# -----------------------------------------------------------------
# 3) Data row 2: regular Arial/10/black, no fill,
#    thin borders - dark top (touches header), light elsewhere
# -----------------------------------------------------------------
for ($col = 1; $col -le 8; $col++) {
  $cell = $ws.Cells.Item(2, $col)
  $cell.Font.Name = "Arial"
  $cell.Font.Size = 10
  $cell.Font.Bold = $false
  $cell.Font.Color = $cBlack
  $cell.Interior.Pattern = -4142
  $cell.NumberFormat = "General"
  $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeLeft).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeTop).Color = $cBorderDk
  $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeRight).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeBottom).Color = $cBorderLt
}

# -----------------------------------------------------------------
# 4) Data row 3: regular Arial/10/black, no fill, thin light borders
# -----------------------------------------------------------------
for ($col = 1; $col -le 8; $col++) {
  $cell = $ws.Cells.Item(3, $col)
  $cell.Font.Name = "Arial"
  $cell.Font.Size = 10
  $cell.Font.Bold = $false
  $cell.Font.Color = $cBlack
  $cell.Interior.Pattern = -4142
  $cell.NumberFormat = "General"
  $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeLeft).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeTop).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeRight).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeBottom).Color = $cBorderLt
}

# -----------------------------------------------------------------
# 5) Row 4 (blank spacer row): default Arial/10/theme font,
#    thin borders - dark/light top (touches row 3), light sides/bottom
# -----------------------------------------------------------------
for ($col = 1; $col -le 8; $col++) {
  $cell = $ws.Cells.Item(4, $col)
  $cell.Font.Name = "Arial"
  $cell.Font.Size = 10
  $cell.Font.Bold = $false
  $cell.Interior.Pattern = -4142
  $cell.NumberFormat = "General"
  $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeLeft).Color = $cBorderLt2
  $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeTop).Color = $cBorderLt
  $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeRight).Color = $cBorderLt2
  $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
  $cell.Borders.Item($xlEdgeBottom).Color = $cBorderLt2
}

# -----------------------------------------------------------------
# 6) Rows 5-10 (blank spacer rows): same font, thin light-grey borders
# -----------------------------------------------------------------
for ($row = 5; $row -le 10; $row++) {
  for ($col = 1; $col -le 8; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Bold = $false
    $cell.Interior.Pattern = -4142
    $cell.NumberFormat = "General"
    $cell.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeLeft).Color = $cBorderLt2
    $cell.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeTop).Color = $cBorderLt2
    $cell.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeRight).Color = $cBorderLt2
    $cell.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $cell.Borders.Item($xlEdgeBottom).Color = $cBorderLt2
  }
}

# -----------------------------------------------------------------
# 7) Column widths (slightly adjusted per the template update)
# -----------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.38
$ws.Columns.Item(2).ColumnWidth = 30.5
$ws.Columns.Item(3).ColumnWidth = 25.63
$ws.Columns.Item(4).ColumnWidth = 23.75
$ws.Columns.Item(5).ColumnWidth = 18.88
$ws.Columns.Item(6).ColumnWidth = 20.75
$ws.Columns.Item(7).ColumnWidth = 21.88
$ws.Columns.Item(8).ColumnWidth = 21.25
